# Sort the comma-separated "Recorded By" names in column G alphabetically
# (case-insensitive, stable), for every data row of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $parts = $text -split ',\s*' | ForEach-Object { $_.Trim() }
    if ($parts.Count -le 1) {
        continue
    }

    $sortedParts = $parts | Sort-Object
    $newText = $sortedParts -join ', '

    if ($newText -ne $text) {
        $cell.Value = $newText
    }
}
